# Apply the edits described by the diff:
# 1. Slide 2: "Summary of Procedures" -> "Summary of Procedure"
# 2. Slide 3: "Handle ECMP for SR Paths" -> "Handle ECMP for SR paths"
# 3. Slide 4: remove the "Similar to the widely deployed synthetic packet loss
#    metric" bullet, and simplify the "Use test packet term for test packet
#    packet, Sender as Session-Sender" bullet text.

$p = $ppt.ActivePresentation

# --- Slide 2: "Summary of Procedures" -> "Summary of Procedure" ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$para = $tr2.Paragraphs(2, 1)
$para.Text = "Summary of Procedure"

# --- Slide 3: "Handle ECMP for SR Paths" -> "Handle ECMP for SR paths" ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$para = $tr3.Paragraphs(6, 1)
$para.Text = "Handle ECMP for SR paths"

# --- Slide 4: remove a bullet and tweak the text of another ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Remove paragraph 7: "Similar to the widely deployed synthetic packet loss metric"
$para = $tr4.Paragraphs(7, 1)
$para.Delete()

# After deletion, the old paragraph 10 ("Use test packet term for test packet
# packet, Sender as Session-Sender") is now paragraph 9.
$para = $tr4.Paragraphs(9, 1)
$para.Text = "Use test packet term, Sender as Session-Sender"
